$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: cells below assigned numeric-looking text values are forced to Text format
# so Excel keeps them as strings (matching source data), then style is reset to
# "Normal" so no stray number-format style index is left on the cell.

# Row 43 and 44 special-case: Kaspa/VeChain rows swapped with new values
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0515"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +15.41%  "

$ws.Range("B44").Value = "Kaspa"
$ws.Range("C44").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.163"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.99%  "

# Updated price / volume values
$ws.Range("D2").Value = "96.689.26"
$ws.Range("E2").Value = "  +0.19%  "
$ws.Range("D3").Value = "3.687.62"
$ws.Range("E3").Value = "  -0.41%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "2.42"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +28.06%  "
$ws.Range("E5").Value = "  +0.03%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "228.87"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.35%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "651.87"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  +2.30%  "
$ws.Range("E9").Value = "  +8.29%  "
$ws.Range("D11").Value = "3.684.95"
$ws.Range("E11").Value = "  -0.45%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "48.04"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +8.15%  "
$ws.Range("E13").Value = "  +2.42%  "
$ws.Range("E14").Value = "  -1.44%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.59"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.72%  "
$ws.Range("D16").Value = "4.377.79"
$ws.Range("E16").Value = "  -0.28%  "
$ws.Range("D17").Value = "96.362.55"
$ws.Range("E17").Value = "  +0.11%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "8.88"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.61%  "
$ws.Range("D19").Value = "3.681.85"
$ws.Range("E19").Value = "  -0.65%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.15"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.46%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.87"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.97%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.545"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "528.65"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.99%  "
$ws.Range("E24").Value = "  -2.00%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.244"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +41.55%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "119.17"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +18.15%  "
$ws.Range("E27").Value = "  +0.88%  "
$ws.Range("E28").Value = "  -2.28%  "
$ws.Range("D29").Value = "3.885.60"
$ws.Range("E29").Value = "  -0.36%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "12.90"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.78%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "13.32"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +9.88%  "
$ws.Range("E32").Value = "  -0.77%  "
$ws.Range("E33").Value = "  +0.00%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.186"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.88%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "33.23"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.20%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.81"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.20%  "
$ws.Range("E37").Value = "  -0.14%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.612"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.28%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "610.82"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -6.77%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.42"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.84%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.09"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.99%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.491"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +13.45%  "
$ws.Range("E46").Value = "  -1.92%  "
$ws.Range("E47").Value = "  -3.87%  "
$ws.Range("E48").Value = "  -0.11%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.97"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +6.08%  "
$ws.Range("E50").Value = "  +0.00%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "23.54"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.14%  "
